$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1944444444444444
$ws.Range("C2").Value = 0.5401234567901234
$ws.Range("J2").Value = 0.02160493827160494
$ws.Range("P2").Value = 0.1481481481481481
$ws.Range("S2").Value = 0.09567901234567901
$ws.Range("B3").Value = 0.02512562814070352
$ws.Range("C3").Value = 0.04020100502512563
$ws.Range("J3").Value = 0.02512562814070352
$ws.Range("P3").Value = 0.7437185929648241
$ws.Range("S3").Value = 0.1658291457286432
$ws.Range("J4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.8484848484848485
$ws.Range("S4").Value = 0.1212121212121212
$ws.Range("P5").Value = 0.8333333333333334
$ws.Range("S5").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.03603603603603604
$ws.Range("F6").Value = 0.04504504504504504
$ws.Range("J6").Value = 0.2432432432432433
$ws.Range("O6").Value = 0.04504504504504504
$ws.Range("Q6").Value = 0.1486486486486487
$ws.Range("R6").Value = 0.07207207207207207
$ws.Range("S6").Value = 0.4099099099099099
$ws.Range("B7").Value = 0.1105527638190955
$ws.Range("F7").Value = 0.02512562814070352
$ws.Range("J7").Value = 0.1155778894472362
$ws.Range("O7").Value = 0.02010050251256281
$ws.Range("Q7").Value = 0.1809045226130653
$ws.Range("R7").Value = 0.07537688442211055
$ws.Range("S7").Value = 0.4723618090452261
$ws.Range("B8").Value = 0.09927360774818401
$ws.Range("D8").Value = 0.01452784503631961
$ws.Range("F8").Value = 0.05326876513317191
$ws.Range("J8").Value = 0.1234866828087167
$ws.Range("O8").Value = 0.007263922518159807
$ws.Range("Q8").Value = 0.1162227602905569
$ws.Range("R8").Value = 0.1089588377723971
$ws.Range("S8").Value = 0.4769975786924939
$ws.Range("B9").Value = 0.1137724550898204
$ws.Range("D9").Value = 0.01796407185628742
$ws.Range("E9").Value = 0.005988023952095809
$ws.Range("F9").Value = 0.08982035928143713
$ws.Range("J9").Value = 0.1137724550898204
$ws.Range("O9").Value = 0.005988023952095809
$ws.Range("Q9").Value = 0.125748502994012
$ws.Range("R9").Value = 0.0658682634730539
$ws.Range("S9").Value = 0.4610778443113773
$ws.Range("B10").Value = 0.1258964143426295
$ws.Range("D10").Value = 0.01832669322709163
$ws.Range("E10").Value = 0.00398406374501992
$ws.Range("F10").Value = 0.08207171314741035
$ws.Range("J10").Value = 0.1115537848605578
$ws.Range("O10").Value = 0.01195219123505976
$ws.Range("Q10").Value = 0.1816733067729084
$ws.Range("R10").Value = 0.07250996015936255
$ws.Range("S10").Value = 0.3920318725099601
$ws.Range("G11").Value = 0.1545189504373178
$ws.Range("J11").Value = 0.08746355685131195
$ws.Range("K11").Value = 0.2011661807580175
$ws.Range("L11").Value = 0.5276967930029155
$ws.Range("S11").Value = 0.02915451895043732
$ws.Range("G12").Value = 0.723404255319149
$ws.Range("J12").Value = 0.1914893617021277
$ws.Range("L12").Value = 0.02659574468085106
$ws.Range("S12").Value = 0.05851063829787234
$ws.Range("G13").Value = 0.5853658536585366
$ws.Range("J13").Value = 0.3414634146341464
$ws.Range("S13").Value = 0.07317073170731707
$ws.Range("F15").Value = 0.03448275862068965
$ws.Range("H15").Value = 0.1231527093596059
$ws.Range("I15").Value = 0.0541871921182266
$ws.Range("J15").Value = 0.3596059113300493
$ws.Range("K15").Value = 0.1231527093596059
$ws.Range("M15").Value = 0.009852216748768473
$ws.Range("O15").Value = 0.04926108374384237
$ws.Range("S15").Value = 0.2463054187192118
$ws.Range("F16").Value = 0.02304147465437788
$ws.Range("H16").Value = 0.1474654377880184
$ws.Range("I16").Value = 0.05990783410138249
$ws.Range("J16").Value = 0.3778801843317972
$ws.Range("K16").Value = 0.1382488479262673
$ws.Range("M16").Value = 0.02764976958525346
$ws.Range("N16").Value = 0.004608294930875576
$ws.Range("O16").Value = 0.05529953917050692
$ws.Range("S16").Value = 0.1658986175115207
$ws.Range("F17").Value = 0.0217983651226158
$ws.Range("H17").Value = 0.1416893732970027
$ws.Range("I17").Value = 0.08446866485013624
$ws.Range("J17").Value = 0.4686648501362398
$ws.Range("K17").Value = 0.08446866485013624
$ws.Range("M17").Value = 0.01362397820163488
$ws.Range("O17").Value = 0.06539509536784741
$ws.Range("S17").Value = 0.1198910081743869
$ws.Range("H18").Value = 0.1783783783783784
$ws.Range("I18").Value = 0.05405405405405406
$ws.Range("J18").Value = 0.4540540540540541
$ws.Range("K18").Value = 0.07027027027027027
$ws.Range("M18").Value = 0.02702702702702703
$ws.Range("O18").Value = 0.03243243243243243
$ws.Range("S18").Value = 0.1837837837837838
$ws.Range("F19").Value = 0.0178173719376392
$ws.Range("H19").Value = 0.2048997772828508
$ws.Range("I19").Value = 0.07720861172976985
$ws.Range("J19").Value = 0.3608017817371937
$ws.Range("K19").Value = 0.1328878990348923
$ws.Range("M19").Value = 0.01707498144023756
$ws.Range("N19").Value = 0.001484780994803266
$ws.Range("O19").Value = 0.066815144766147
$ws.Range("S19").Value = 0.1210096510764662
